$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column C (the "Förändrad" / Changed date column)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

# Update every populated cell in column C (rows 2..lastRow) from 45179 to 45180
$range = $ws.Range("C2:C$lastRow")
for ($i = 1; $i -le $range.Rows.Count; $i++) {
    $cell = $range.Cells.Item($i, 1)
    if ($cell.Value2 -eq 45179) {
        $cell.Value2 = 45180
    }
}

$wb.Save()
